$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = "Trophic level"
$ws.Range("E10").Value = "Trophic level will automatically be put on taxon level if size classes are equal. Are different for Unicell etc."

$ws.Range("D46").Select()
